# Adds the three new character styles (GaNStyle, GaNParagraph, GaNLinks)
# and applies them to the relevant runs, per the commit's intent of
# "Add styles to the new paragraphs".

$d = $word.ActiveDocument

# --- Define the new character styles -------------------------------------

$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Dates à utiliser..." run (4 occurrences) ---

$datesText = "Dates à utiliser pour la Campagne 2022 Constellation des Gémeaux: 14-23 février, 14-24 mars"
$rng = $d.Content
$rng.Find.ClearFormatting()
while ($rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- Apply GaNParagraph to the campaign description paragraph ------------

$paragraphText = "Vous allez participer à une campagne mondiale d’observation pour détecter les plus faibles étoiles visibles afin de mesurer la pollution lumineuse sur un site donné. Partout dans le monde, en localisant et en observant la Constellation des Gémeaux dans le ciel nocturne et en la comparant aux cartes stellaires, les participants, apprendront comment l’éclairage, dans leur environnement local, influence la pollution lumineuse. Vos contributions à la base de données en ligne permettront de mesurer la qualité du ciel nocturne."
$rng = $d.Content
$rng.Find.ClearFormatting()
if ($rng.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the map credit line --------------------------------

$linksText = "Les cartes figurant dans ce document ont été établies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng = $d.Content
$rng.Find.ClearFormatting()
if ($rng.Find.Execute($linksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
}
